$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.1635273785303319
$ws.Range("J2").Value = 0.1635273785303318
$ws.Range("O2").Value = 0.4715772180415435
$ws.Range("P2").Value = 0.4715772180415434
$ws.Range("S2").Value = 0.07711578624096033
$ws.Range("T2").Value = 0.0771157862409603

# Row 3
$ws.Range("I3").Value = 0.1635273785303319
$ws.Range("J3").Value = 0.1635273785303318
$ws.Range("M3").Value = 0.2448813333333333
$ws.Range("N3").Value = 0.734644
$ws.Range("O3").Value = 0.5284227819584566
$ws.Range("P3").Value = 0.5284227819584566
$ws.Range("Q3").Value = 0.0006101626555555555
$ws.Range("R3").Value = 0.005491463899999999
$ws.Range("S3").Value = 0.08641159228937155
$ws.Range("T3").Value = 0.08641159228937154

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.01274533333333334
$ws.Range("H4").Value = 0.03823600000000001
$ws.Range("I4").Value = 0.8364726214696682
$ws.Range("J4").Value = 0.8364726214696681
$ws.Range("O4").Value = 0.4715772180415435
$ws.Range("P4").Value = 0.4715772180415434
$ws.Range("Q4").Value = 0.002785339656000001
$ws.Range("R4").Value = 0.02506805690400001
$ws.Range("S4").Value = 0.3944614318005831
$ws.Range("T4").Value = 0.3944614318005831

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.01274533333333334
$ws.Range("H5").Value = 0.03823600000000001
$ws.Range("I5").Value = 0.8364726214696682
$ws.Range("J5").Value = 0.8364726214696681
$ws.Range("M5").Value = 0.2448813333333333
$ws.Range("N5").Value = 0.734644
$ws.Range("O5").Value = 0.5284227819584566
$ws.Range("P5").Value = 0.5284227819584566
$ws.Range("Q5").Value = 0.003121094220444445
$ws.Range("R5").Value = 0.028089847984
$ws.Range("S5").Value = 0.442011189669085
$ws.Range("T5").Value = 0.442011189669085
